# Auto-generated edit script: apply updated '想去人数' (want-to-go count) values
# and one venue-address correction, across all four worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 27
$ws.Cells.Item(3, 6).Value = 788
$ws.Cells.Item(7, 6).Value = 1148
$ws.Cells.Item(8, 6).Value = 903
$ws.Cells.Item(10, 6).Value = 717
$ws.Cells.Item(11, 6).Value = 1032
$ws.Cells.Item(12, 6).Value = 1430
$ws.Cells.Item(15, 6).Value = 1599
$ws.Cells.Item(17, 6).Value = 605
$ws.Cells.Item(23, 6).Value = 749
$ws.Cells.Item(24, 6).Value = 610
$ws.Cells.Item(25, 6).Value = 486
$ws.Cells.Item(27, 4).Value = "长宁路1191号来福士西区(W)B1层01号、11号 星零界"
$ws.Cells.Item(30, 6).Value = 298
$ws.Cells.Item(31, 6).Value = 2416
$ws.Cells.Item(32, 6).Value = 278
$ws.Cells.Item(34, 6).Value = 459
$ws.Cells.Item(36, 6).Value = 3941
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(11, 6).Value = 8
$ws.Cells.Item(12, 6).Value = 395
$ws.Cells.Item(20, 6).Value = 255
$ws.Cells.Item(22, 6).Value = 120
$ws.Cells.Item(28, 6).Value = 1711
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 95
$ws.Cells.Item(4, 6).Value = 1268
$ws.Cells.Item(5, 6).Value = 1660
$ws.Cells.Item(7, 6).Value = 996
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 95
$ws.Cells.Item(3, 6).Value = 1268
$ws.Cells.Item(4, 6).Value = 1660
$ws.Cells.Item(6, 6).Value = 996
$ws.Cells.Item(7, 6).Value = 27
$ws.Cells.Item(8, 6).Value = 788
$ws.Cells.Item(12, 6).Value = 1148
$ws.Cells.Item(13, 6).Value = 903
$ws.Cells.Item(17, 6).Value = 717
$ws.Cells.Item(21, 6).Value = 1032
$ws.Cells.Item(22, 6).Value = 1430
$ws.Cells.Item(25, 6).Value = 1599
$ws.Cells.Item(27, 6).Value = 605
$ws.Cells.Item(32, 6).Value = 749
$ws.Cells.Item(33, 6).Value = 610
$ws.Cells.Item(34, 6).Value = 486
$ws.Cells.Item(36, 4).Value = "长宁路1191号来福士西区(W)B1层01号、11号 星零界"
$ws.Cells.Item(38, 6).Value = 255
$ws.Cells.Item(42, 6).Value = 298
$ws.Cells.Item(43, 6).Value = 2416
$ws.Cells.Item(46, 6).Value = 1711
$ws.Cells.Item(47, 6).Value = 1711
$ws.Cells.Item(49, 6).Value = 459
$ws.Cells.Item(50, 6).Value = 3941
